# Generate Report for Handback
#
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   (Overview + per-language sheets)
# - The "Latest Handback DateTime" placeholder (0001-01-01 00:00:00) is
#   replaced with the real handback timestamp for each language
# - The newly-populated "Latest Target File" / "Latest Handback File"
#   columns (F/G) are filled in with hyperlinked file names now that the
#   handback package exists

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$mdFile  = "819d8295-52b8-4b3c-aaf7-ffdb05d1d164.md"
$xlfFileZh = "819d8295-52b8-4b3c-aaf7-ffdb05d1d164.dc06b4015d182ad9f6866184d668cd1729fb8b1b.zh-cn.xlf"
$xlfFileDe = "819d8295-52b8-4b3c-aaf7-ffdb05d1d164.dc06b4015d182ad9f6866184d668cd1729fb8b1b.de-de.xlf"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/49275fa2d37c31f6651fbdf9e9159c66827769f1/e2e/819d8295-52b8-4b3c-aaf7-ffdb05d1d164.md"
$xlfUrlZh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/004f14cca44dc33499c68618d0dbf52b04e60a5c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/819d8295-52b8-4b3c-aaf7-ffdb05d1d164.dc06b4015d182ad9f6866184d668cd1729fb8b1b.zh-cn.xlf"
$xlfUrlDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1df2dc39c95873e899a096284539068c1a8bef5d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/819d8295-52b8-4b3c-aaf7-ffdb05d1d164.dc06b4015d182ad9f6866184d668cd1729fb8b1b.de-de.xlf"

# ---- Overview sheet: flip the status shown for both languages ----
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($row in 2, 3) {
    foreach ($col in "B", "C") {
        $cell = $wsOverview.Range($col + $row)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($row in 2, 3) {
    $wsZh.Range("C" + $row).Value = $newStatus
    $wsZh.Range("H" + $row).Value = "2016-03-18 05:59:08"
    $wsZh.Hyperlinks.Add($wsZh.Range("F" + $row), $mdUrl, "", "", $mdFile) | Out-Null
    $wsZh.Hyperlinks.Add($wsZh.Range("G" + $row), $xlfUrlZh, "", "", $xlfFileZh) | Out-Null
}

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($row in 2, 3) {
    $wsDe.Range("C" + $row).Value = $newStatus
    $wsDe.Range("H" + $row).Value = "2016-03-18 05:59:13"
    $wsDe.Hyperlinks.Add($wsDe.Range("F" + $row), $mdUrl, "", "", $mdFile) | Out-Null
    $wsDe.Hyperlinks.Add($wsDe.Range("G" + $row), $xlfUrlDe, "", "", $xlfFileDe) | Out-Null
}
